# "checkout one scenario completed"
# Adds country/city/address/zipCode/mob columns (D:H) with a second test
# record, widens column C and H, and moves the active selection to H3
# (mirroring where the user's cursor ended up after typing the new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new header row (row 1), columns D:H -----------------------------
$ws.Range("D1").Value = "country"
$ws.Range("E1").Value = "city"
$ws.Range("F1").Value = "address"
$ws.Range("G1").Value = "zipCode"
$ws.Range("H1").Value = "mob"

# --- new data row (row 2), columns D:H --------------------------------
$ws.Range("D2").Value = "India"
$ws.Range("E2").Value = "Bengaluru"
$ws.Range("F2").Value = "Bengaluru"

# zip code / mobile number are numeric-looking text, so they need the
# same quote-prefixed text style the existing "password1" cell (B2) uses
$ws.Range("G2").Value = "'560001"
$ws.Range("H2").Value = "'9999999999"

# --- column sizing -----------------------------------------------------
$ws.Columns("C").ColumnWidth = 31.33
$ws.Columns("H").ColumnWidth = 10.166666666666666

# --- selection mirrors where the user left the cursor -------------------
$ws.Range("H3").Select() | Out-Null
